$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as plain text in the source data (e.g. "609.14"),
# using "." as both a thousands separator and decimal point. Force text storage so
# Excel does not reinterpret these as numbers, then restore the default "Normal" style
# so no stray cell formatting is introduced.
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D26","D27","D28","D29","D32","D33","D34","D37","D38","D42","D43","D44","D46","D47","D48","D49","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.202.29"
$ws.Range("E2").Value = "  -3.67%  "
$ws.Range("D3").Value = "3.491.53"
$ws.Range("E3").Value = "  -5.21%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "609.02"
$ws.Range("E5").Value = "  -6.49%  "
$ws.Range("D6").Value = "148.75"
$ws.Range("E6").Value = "  -7.99%  "
$ws.Range("D7").Value = "3.490.60"
$ws.Range("E7").Value = "  -5.17%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  -3.43%  "
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -5.34%  "
$ws.Range("D11").Value = "6.93"
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("E13").Value = "  -7.21%  "
$ws.Range("D14").Value = "4.079.98"
$ws.Range("E14").Value = "  -5.24%  "
$ws.Range("D15").Value = "31.47"
$ws.Range("E15").Value = "  -3.98%  "
$ws.Range("D16").Value = "3.500.68"
$ws.Range("E16").Value = "  -4.66%  "
$ws.Range("D17").Value = "67.098.36"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "6.45"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").Value = "15.03"
$ws.Range("E20").Value = "  -5.90%  "
$ws.Range("D21").Value = "446.69"
$ws.Range("E21").Value = "  -5.26%  "
$ws.Range("D22").Value = "9.03"
$ws.Range("E22").Value = "  -12.72%  "
$ws.Range("D23").Value = "0.624"
$ws.Range("E23").Value = "  -5.39%  "
$ws.Range("D24").Value = "77.13"
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "3.630.74"
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  -9.26%  "
$ws.Range("D29").Value = "8.26"
$ws.Range("E29").Value = "  -6.92%  "
$ws.Range("E30").Value = "  -4.61%  "
$ws.Range("E31").Value = "  -7.68%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "0.165"
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("D34").Value = "25.70"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("E36").Value = "  -7.40%  "
$ws.Range("D37").Value = "3.476.10"
$ws.Range("E37").Value = "  -5.51%  "
$ws.Range("D38").Value = "8.00"
$ws.Range("E38").Value = "  -4.92%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "0.0871"
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("D43").Value = "168.32"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").Value = "5.44"
$ws.Range("E44").Value = "  -7.78%  "
$ws.Range("E45").Value = "  -5.10%  "
$ws.Range("D46").Value = "45.42"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").Value = "1.25"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.54"
$ws.Range("E48").Value = "  -9.96%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "26.28"
$ws.Range("E49").Value = "  -9.95%  "
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  -4.15%  "

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
